# OLX Monitor 2026-02-17 14:21 update
# Appends a fresh monitoring run (8 listing rows) to the "PODSUMOWANIE" log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

$timestamp = "2026-02-17 14:21:01"

# New rows of scraped listing data for this run.
$rows = @(
    @{ B = "poqui";           C = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda";               D = 2299;  E = "19.01.2026"; F = 29;  G = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html";                               H = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR" },
    @{ B = "poqui";           C = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";          D = 2499;  E = "28.10.2025"; F = 111; G = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html";                           H = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" },
    @{ B = "poqui";           C = "Przytulny pokój blisko Politechniki – ul. Przytulna";                     D = 599;   E = "10.10.2025"; F = 130; G = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html";                                     H = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" },
    @{ B = "poqui";           C = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza";           D = 2049;  E = "19.12.2025"; F = 59;  G = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html";                            H = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc" },
    @{ B = "pokojewlublinie"; C = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                     D = 58640; E = "11.08.2025"; F = 190; G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html";                                 H = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" },
    @{ B = "pokojewlublinie"; C = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";   D = 12640; E = "19.01.2026"; F = 29;  G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html";                  H = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" },
    @{ B = "dawnypatron";     C = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";   D = 730;   E = "20.09.2024"; F = 515; G = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html";              H = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" },
    @{ B = "dawnypatron";     C = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";  D = 14690; E = "05.12.2025"; F = 74;  G = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"; H = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

$startRow = 39
$endRow = $startRow + $rows.Count - 1

# Copy formatting from the previous run's block (rows 31-38) down onto the
# freshly appended rows (39-46) so styles (borders/alignment/highlighting)
# match the established pattern.
$srcFormat = $ws.Range("A31:H38")
$dstFormat = $ws.Range("A$startRow`:H$endRow")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    # Column E holds a "dd.mm.yyyy" looking string that must stay plain text
    # (as in the source data). Excel's smart entry would silently reinterpret
    # unambiguous values (day and month both <=12) as real dates, so force
    # text formatting while assigning it.
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row.E

    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Re-apply the reference formatting once more. Forcing the "@" text format on
# column E above can leave a stray quote-prefix style on those cells; pasting
# the original formats back over the whole block restores the exact style
# indices (13/14/15) used by the rest of the table without touching values.
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
